$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B48 was stored as text "3"; convert it to a real number 3.
$ws.Range("B48").Value = 3

# Add new row 49 with the additional annotation record.
$ws.Range("A49").Value = "Ying Tang"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "1"
$ws.Range("B49").ClearFormats()
$ws.Range("C49").Value = "It is incredible that,so rude and misleading,willful misinterpretations and falsehoods "
$ws.Range("D49").Value = "CRT"
$ws.Range("E49").Value = "OTH"
$ws.Range("F49").Value = "f6e31c12-680e-4edf-b736-d4a426f6f32f"
$ws.Range("G49").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H49").Value = "It is incredible that the commenter continues to be so rude and misleading (should OpenReview have a moderating system?), and continues to frame this interaction as an attempt to convince *them* rather than to correct the constant series of willful misinterpretations and falsehoods that they manage to state about our work in every single interaction, in the hope that they do not mislead others."
